# calibration_prot.xlsx - update calibration data points and refresh selection
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 (A10=2.4): recalculated X600 / X600_cor values
$ws.Cells.Item(10, 2).Value = 0.505
$ws.Cells.Item(10, 3).Value = 0.45450000000000002

# Row 11 (A11=3): recalculated X600 / X600_cor values
$ws.Cells.Item(11, 2).Value = 0.55500000000000005
$ws.Cells.Item(11, 3).Value = 0.50450000000000006

# Leave the active cell where the author last left it when saving
$ws.Range("E10").Select()
